# Update the document title from "TRADUCKXION V2.5" to "TRADUCKXION V2.6".
#
# The target OOXML keeps the unchanged prefix ("...V2.") in the original
# run and puts the new "6" in a brand-new run (both runs end up with the
# same <w:rPr><w:lang w:val="fr-FR"/></w:rPr>). A plain text/Find-Replace
# assignment on the digit would just rewrite the existing run in place
# (and the engine re-merges runs that end up with identical formatting),
# so we briefly flip a character-formatting property on the replacement
# text to force it into its own run, then clear that formatting again so
# the final run properties match the original ("lang fr-FR" only).

$d = $word.ActiveDocument

$titleSearch = $d.Content
$found = $titleSearch.Find.Execute("TRADUCKXION V2.5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Narrow the match down to just the final character, the "5" in "V2.5".
    $digit = $d.Range($titleSearch.End - 1, $titleSearch.End)

    $digit.Bold = 1
    $digit.Text = "6"

    # Re-fetch the range for the freshly typed "6" (Range positions stay
    # valid since the replacement text is the same length) and drop the
    # temporary bold flag so it ends up with only the inherited language
    # formatting, same as the run it was split from.
    $newDigit = $d.Range($titleSearch.End - 1, $titleSearch.End)
    $newDigit.Bold = 0
}
